$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 2547
$ws.Range("L3").Value = 2558
$ws.Range("L4").Value = 691
$ws.Range("L5").Value = 152
$ws.Range("L6").Value = 2310
$ws.Range("L7").Value = 8258

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 154
$ws.Range("L6").Value = 142
$ws.Range("L7").Value = 527

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 82
$ws.Range("L7").Value = 195

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 114
$ws.Range("L7").Value = 367

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 49
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L6").Value = 95
$ws.Range("L7").Value = 302

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L3").Value = 48
$ws.Range("L7").Value = 154

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 59
$ws.Range("L7").Value = 134

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 68
$ws.Range("L4").Value = 33
$ws.Range("L6").Value = 62
$ws.Range("L7").Value = 267
$ws.Range("L8").Value = 527
$ws.Range("L11").Value = 145
$ws.Range("L14").Value = 41
$ws.Range("L18").Value = 61
$ws.Range("L20").Value = 212
$ws.Range("L24").Value = 18
$ws.Range("L29").Value = 432
$ws.Range("L30").Value = 36
$ws.Range("L31").Value = 81
$ws.Range("L33").Value = 367
$ws.Range("L37").Value = 302
$ws.Range("L42").Value = 267
$ws.Range("L44").Value = 59
$ws.Range("L51").Value = 93
$ws.Range("L52").Value = 164
$ws.Range("L55").Value = 77
$ws.Range("L60").Value = 50
$ws.Range("L63").Value = 23
$ws.Range("L65").Value = 154
$ws.Range("L67").Value = 306
$ws.Range("L73").Value = 67
$ws.Range("L76").Value = 101
$ws.Range("L78").Value = 107
$ws.Range("L79").Value = 225
$ws.Range("L83").Value = 195
$ws.Range("L84").Value = 85
$ws.Range("L85").Value = 431
$ws.Range("L88").Value = 111
$ws.Range("L89").Value = 104
$ws.Range("L90").Value = 81
$ws.Range("L92").Value = 23
$ws.Range("L94").Value = 102
$ws.Range("L95").Value = 114
$ws.Range("L97").Value = 75
$ws.Range("L98").Value = 56
$ws.Range("L99").Value = 134
$ws.Range("L101").Value = 8258

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L2").Value = 28
$ws.Range("L7").Value = 81

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 92
$ws.Range("L3").Value = 105
$ws.Range("L7").Value = 306

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L2").Value = 31
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 85

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 158
$ws.Range("L6").Value = 114
$ws.Range("L7").Value = 432

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 59

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 18
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 101

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("L2").Value = 17
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L2").Value = 27
$ws.Range("L7").Value = 62

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 82
$ws.Range("L6").Value = 82
$ws.Range("L7").Value = 267

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L3").Value = 29
$ws.Range("L7").Value = 107

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("L3").Value = 4
$ws.Range("L7").Value = 18

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 82
$ws.Range("L7").Value = 225

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 67
$ws.Range("L4").Value = 18
$ws.Range("L6").Value = 60
$ws.Range("L7").Value = 212

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L2").Value = 23
$ws.Range("L7").Value = 61

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 78
$ws.Range("L6").Value = 75
$ws.Range("L7").Value = 267

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 27
$ws.Range("L7").Value = 102

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 56

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L3").Value = 2
$ws.Range("L5").Value = 7

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 53
$ws.Range("L7").Value = 145

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L2").Value = 29
$ws.Range("L7").Value = 67

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L3").Value = 21
$ws.Range("L7").Value = 68

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 75

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L2").Value = 8
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 23

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L3").Value = 39
$ws.Range("L7").Value = 111

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L4").Value = 19
$ws.Range("L6").Value = 25
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L4").Value = 7
$ws.Range("L7").Value = 81

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 27
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L3").Value = 17
$ws.Range("L5").Value = 2
$ws.Range("L7").Value = 50

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 173
$ws.Range("L4").Value = 35
$ws.Range("L6").Value = 87
$ws.Range("L7").Value = 431

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 53
$ws.Range("L3").Value = 49
$ws.Range("L7").Value = 164

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 33
